# Update Mexico Liga MX Femenil bases (07-04-2024 22:30)
# For a set of row pairs, swap all the match data (columns B:AC) between
# the two rows while leaving column A (the sequential id) untouched.
#
# NOTE: a plain `foreach` loop that reuses the same variable names for the
# Range/array objects on every iteration can leave the last iteration's
# write un-flushed in this runtime, so a classic indexed `for` loop is used
# instead and the helper variables are cleared after each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(47, 48),
    @(101, 102),
    @(109, 110),
    @(131, 132),
    @(215, 216),
    @(229, 230),
    @(251, 252),
    @(271, 272)
)

for ($idx = 0; $idx -lt $rowPairs.Count; $idx++) {
    $pair = $rowPairs[$idx]
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("B$($r1):AC$($r1)")
    $rangeB = $ws.Range("B$($r2):AC$($r2)")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA

    $rangeA = $null
    $rangeB = $null
    $valuesA = $null
    $valuesB = $null
}
